# Apply the LOM3208.xlsx content update.
# The original row 13 (blank A cell, "144651 - Antonio Fernando Sartori" in B/C)
# is removed entirely, shifting all subsequent rows up by one. Several cells
# further down the sheet are then updated with new text content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old row 13 ("Docentes responsáveis:" value row) - this shifts
# every row below it up by one, matching the target layout/row-heights.
$ws.Rows(13).Delete()

# After the shift, update the cells whose text content changed.
$ws.Range("B10").Value = "144651 - Antonio Fernando Sartori"
$ws.Range("C10").Value = "144651 - Antonio Fernando Sartori"

$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

$ws.Range("B15").Value = "01/01/2012"
$ws.Range("C15").Value = "01/01/2012"

$ws.Range("B18").Value = "144651 - Antonio Fernando Sartori"
$ws.Range("C18").Value = "144651 - Antonio Fernando Sartori"

$ws.Range("B19").Value = "Aulas expositivas, seminários e exercícios comentados."
$ws.Range("C19").Value = "Aulas expositivas, seminários e exercícios comentados."

$ws.Range("B20").Value = "Média aritmética de duas provas sendo a primeira com peso 1 e a segunda com peso 2."
$ws.Range("C20").Value = "Média aritmética de duas provas sendo a primeira com peso 1 e a segunda com peso 2."

$ws.Range("B21").Value = "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"
$ws.Range("C21").Value = "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"
